{"js": "// Update the date heading and the multiplication-table answers in place.\n// Positional (row/col) targeting is used throughout instead of text search,\n// because several of the old/new answer strings repeat or collide with each\n// other across cells (e.g. \"857\u00d76=5142\" appears twice, and \"258\u00d75=1290\" is\n// simultaneously an old value in one cell and a new value in another).\n\nconst body = context.document.body;\n\n// 1) Heading paragraph: \"2025-06-28 Saturday\" -> \"2025-06-29 Sunday\"\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst heading = paras.items[0];\nheading.insertText(\"2025-06-29 Sunday\", \"Replace\");\n\n// 2) Table of multiplication answers (5 columns; only every 5th row has text)\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst rowUpdates = [\n  { row: 0, values: [\"853\u00d78=6824\", \"638\u00d77=4466\", \"364\u00d78=2912\", \"726\u00d73=2178\", \"418\u00d76=2508\"] },\n  { row: 4, values: [\"864\u00d74=3456\", \"508\u00d75=2540\", \"258\u00d75=1290\", \"975\u00d75=4875\", \"570\u00d75=2850\"] },\n  { row: 9, values: [\"655\u00d79=5895\", \"615\u00d75=3075\", \"219\u00d78=1752\", \"289\u00d76=1734\", \"777\u00d78=6216\"] },\n  { row: 14, values: [\"684\u00d79=6156\", \"773\u00d78=6184\", \"121\u00d79=1089\", \"756\u00d78=6048\", \"435\u00d72=870\"] },\n  { row: 19, values: [\"809\u00d72=1618\", \"436\u00d74=1744\", \"525\u00d72=1050\", \"977\u00d76=5862\", \"468\u00d74=1872\"] },\n];\n\nfor (const { row, values } of rowUpdates) {\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(row, col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the multiplication-table answers in place.\n# Positional (row/col) targeting is used throughout instead of text search,\n# because several of the old/new answer strings repeat or collide with each\n# other across cells (e.g. \"857\u00d76=5142\" appears twice, and \"258\u00d75=1290\" is\n# simultaneously an old value in one cell and a new value in another).\n\n$d = $word.ActiveDocument\n\n# 1) Heading paragraph: \"2025-06-28 Saturday\" -> \"2025-06-29 Sunday\"\n$d.Paragraphs.Item(1).Range.Text = \"2025-06-29 Sunday\"\n\n# 2) Table of multiplication answers (5 columns; only every 5th row has text)\n$t = $d.Tables.Item(1)\n\n$updates = @(\n  @{ Row = 1; Col = 1; Text = \"853\u00d78=6824\" },\n  @{ Row = 1; Col = 2; Text = \"638\u00d77=4466\" },\n  @{ Row = 1; Col = 3; Text = \"364\u00d78=2912\" },\n  @{ Row = 1; Col = 4; Text = \"726\u00d73=2178\" },\n  @{ Row = 1; Col = 5; Text = \"418\u00d76=2508\" },\n  @{ Row = 5; Col = 1; Text = \"864\u00d74=3456\" },\n  @{ Row = 5; Col = 2; Text = \"508\u00d75=2540\" },\n  @{ Row = 5; Col = 3; Text = \"258\u00d75=1290\" },\n  @{ Row = 5; Col = 4; Text = \"975\u00d75=4875\" },\n  @{ Row = 5; Col = 5; Text = \"570\u00d75=2850\" },\n  @{ Row = 10; Col = 1; Text = \"655\u00d79=5895\" },\n  @{ Row = 10; Col = 2; Text = \"615\u00d75=3075\" },\n  @{ Row = 10; Col = 3; Text = \"219\u00d78=1752\" },\n  @{ Row = 10; Col = 4; Text = \"289\u00d76=1734\" },\n  @{ Row = 10; Col = 5; Text = \"777\u00d78=6216\" },\n  @{ Row = 15; Col = 1; Text = \"684\u00d79=6156\" },\n  @{ Row = 15; Col = 2; Text = \"773\u00d78=6184\" },\n  @{ Row = 15; Col = 3; Text = \"121\u00d79=1089\" },\n  @{ Row = 15; Col = 4; Text = \"756\u00d78=6048\" },\n  @{ Row = 15; Col = 5; Text = \"435\u00d72=870\" },\n  @{ Row = 20; Col = 1; Text = \"809\u00d72=1618\" },\n  @{ Row = 20; Col = 2; Text = \"436\u00d74=1744\" },\n  @{ Row = 20; Col = 3; Text = \"525\u00d72=1050\" },\n  @{ Row = 20; Col = 4; Text = \"977\u00d76=5862\" },\n  @{ Row = 20; Col = 5; Text = \"468\u00d74=1872\" }\n)\n\nforeach ($u in $updates) {\n  $t.Cell($u.Row, $u.Col).Range.Text = $u.Text\n}\n"}
